$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that sat after
# "along with the zip folder" (it simply moves, see Change 2 below).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    [void]$d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Change 2: split the run "Create a folder ... Move the zip into this
# folder. Start up the " right after "Move the zip" and drop a new
# "_GoBack" bookmark at the split point.
# ---------------------------------------------------------------------
$marker = "Move the zip"
$tail = " into this folder. Start up the "
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains($marker + $tail)) {
        $pStart = $p.Range.Start
        $splitPos = $pStart + $t.IndexOf($marker) + $marker.Length

        # Insert the (now empty) bookmark exactly at the split point -
        # this naturally breaks the run in two.
        $bmRange = $d.Range($splitPos, $splitPos)
        [void]$d.Bookmarks.Add("_GoBack", $bmRange)

        # Re-stamp the text of the first half so the run's xml:space
        # handling is recalculated now that it no longer ends in a space.
        $firstHalf = $d.Range($pStart, $splitPos)
        $firstText = $firstHalf.Text
        [void]$firstHalf.Delete()
        [void]$firstHalf.InsertBefore($firstText)
        break
    }
}

# ---------------------------------------------------------------------
# Change 3: replace the empty paragraph that follows the "Gather"
# section (just before the "Assess" heading) with new text.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.StartsWith("Assess")) {
            $r = $p.Range
            $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">Alright then, let’s go ahead. Using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zipfile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> module</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
            [void]$r.InsertXML($xml)
            break
        }
    }
}
